$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 467, shifting all
# subsequent rows (old 467:586) down to (468:587).
$ws.Rows.Item(467).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(467, 1).Value2 = 6
$ws.Cells.Item(467, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(467, 3).Value2 = "Metropolitana"
$ws.Cells.Item(467, 4).Value2 = 44932
$ws.Cells.Item(467, 5).Value2 = 13
$ws.Cells.Item(467, 6).Value2 = 100112039
$ws.Cells.Item(467, 7).Value2 = "Ciboulette"
$ws.Cells.Item(467, 8).Value2 = "Sin especificar"
$ws.Cells.Item(467, 9).Value2 = "Primera"
$ws.Cells.Item(467, 10).Value2 = 770
$ws.Cells.Item(467, 11).Value2 = 900
$ws.Cells.Item(467, 12).Value2 = 1000
$ws.Cells.Item(467, 13).Value2 = 945
$ws.Cells.Item(467, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(467, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(467, 16).Value2 = 315
$ws.Cells.Item(467, 17).Value2 = 3
$ws.Cells.Item(467, 18).Value2 = "Hortaliza"
